# Leave Card update 12/22/2023 10:59 AM
# - Fill in PERIOD dates (column A) and EARNED (column C, 1.25) for the
#   11 months Jan-2023..Nov-2023 (rows 80-90) in Table1 on Sheet1.
# - Fill in the PERIOD dates only (column A) for the future months
#   Dec-2023..Feb-2026 (rows 91-117); EARNED stays blank for those.
# - Column G ("EARNED ") is a calculated table column that mirrors C and
#   recalculates automatically.
# - Leave the workbook positioned on the CONVERTION sheet (matches the
#   saved view state in the edited file).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("CONVERTION")

# Table1[PERIOD] serial dates for rows 80..117 (first of each month,
# Jan-2023 through Feb-2026).
$periodDates = @{
    80  = 44927 # 2023-01-01
    81  = 44958 # 2023-02-01
    82  = 44986 # 2023-03-01
    83  = 45017 # 2023-04-01
    84  = 45047 # 2023-05-01
    85  = 45078 # 2023-06-01
    86  = 45108 # 2023-07-01
    87  = 45139 # 2023-08-01
    88  = 45170 # 2023-09-01
    89  = 45200 # 2023-10-01
    90  = 45231 # 2023-11-01
    91  = 45261 # 2023-12-01
    92  = 45292 # 2024-01-01
    93  = 45323 # 2024-02-01
    94  = 45352 # 2024-03-01
    95  = 45383 # 2024-04-01
    96  = 45413 # 2024-05-01
    97  = 45444 # 2024-06-01
    98  = 45474 # 2024-07-01
    99  = 45505 # 2024-08-01
    100 = 45536 # 2024-09-01
    101 = 45566 # 2024-10-01
    102 = 45597 # 2024-11-01
    103 = 45627 # 2024-12-01
    104 = 45658 # 2025-01-01
    105 = 45689 # 2025-02-01
    106 = 45717 # 2025-03-01
    107 = 45748 # 2025-04-01
    108 = 45778 # 2025-05-01
    109 = 45809 # 2025-06-01
    110 = 45839 # 2025-07-01
    111 = 45870 # 2025-08-01
    112 = 45901 # 2025-09-01
    113 = 45931 # 2025-10-01
    114 = 45962 # 2025-11-01
    115 = 45992 # 2025-12-01
    116 = 46023 # 2026-01-01
    117 = 46054 # 2026-02-01
}

# Rows that also get an EARNED credit of 1.25 posted (Jan-2023..Nov-2023).
$earnedRows = 80..90

foreach ($row in 80..117) {
    $ws1.Cells.Item($row, 1).Value = $periodDates[$row]   # column A = PERIOD
    if ($earnedRows -contains $row) {
        $ws1.Cells.Item($row, 3).Value = 1.25              # column C = EARNED
    }
}

# Recalculate so the BALANCE / mirrored EARNED formulas pick up the new
# postings before the view state (and any cached values) are saved.
$excel.Calculate()

# Restore the saved view: Sheet1 scrolled so the frozen/split bottom pane
# starts at row 70 with E93 selected, then CONVERTION left as the active
# (visible) sheet/tab.
$ws1.Activate() | Out-Null
$excel.ActiveWindow.SplitRow = 69
$ws1.Range("E93").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("G3").Select() | Out-Null
